# Generate Report for Handback
#
# Updates the localization-status report after a handback:
#   - "Ready for handoff" status text becomes "Handed back: in sync with en-US"
#     (Overview!E2:F3, and the Status column C on the zh-cn / de-de sheets,
#     since they all share the same string).
#   - On the zh-cn and de-de detail sheets, each of the two data rows gets its
#     "Latest Target File" (I) and "Latest Handback File" (J) columns filled
#     in with the generated markdown / xlf links, and "Latest Handback
#     DateTime" (K) gets stamped with the generation time.
#   - Columns widened to fit the newly-populated hyperlink text.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b6a4b3e70f174ed85c03424e69a84b9257487a65/e2e/"

$file1Md = "15a85270-f4b6-4ef9-817f-95b00a05f115.md"
$file2Md = "f1e83fcc-942a-4159-8ece-a0eb768612ef.md"

$file1ZhXlf = "15a85270-f4b6-4ef9-817f-95b00a05f115.e394e227fec69c8338d4be34f6b555cd0e680858.zh-cn.xlf"
$file2ZhXlf = "f1e83fcc-942a-4159-8ece-a0eb768612ef.990376d0bb0c6b512c26b733c80047f8ed434045.zh-cn.xlf"
$file1DeXlf = "15a85270-f4b6-4ef9-817f-95b00a05f115.e394e227fec69c8338d4be34f6b555cd0e680858.de-de.xlf"
$file2DeXlf = "f1e83fcc-942a-4159-8ece-a0eb768612ef.990376d0bb0c6b512c26b733c80047f8ed434045.de-de.xlf"

$zhDateTime = "2016-08-28 10:56:12"
$deDateTime = "2016-08-28 10:56:19"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the two status cells per row (E/F columns)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn detail sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = $file1ZhXlf
$wsZh.Range("J3").Value = $file2ZhXlf

$wsZh.Range("K2").Value = $zhDateTime
$wsZh.Range("K3").Value = $zhDateTime

# Re-create hyperlinks so A2/A3 keep their relationship slots and I2/I3 get
# new ones for the "Latest Target File" column, in the same order Excel
# would emit them (A2, I2, A3, I3).
$wsZh.Range("A1").Hyperlinks.Delete()

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), ($repoBase + $file2Md), "", "", $file2Md)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------
# de-de detail sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = $file1DeXlf
$wsDe.Range("J3").Value = $file2DeXlf

$wsDe.Range("K2").Value = $deDateTime
$wsDe.Range("K3").Value = $deDateTime

$wsDe.Range("A1").Hyperlinks.Delete()

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($repoBase + $file1Md), "", "", $file1Md)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), ($repoBase + $file2Md), "", "", $file2Md)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($repoBase + $file2Md), "", "", $file2Md)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Output "done"
